$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'26.909.86"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.85%  '
$ws.Range('D3').Value = "'1.862.80"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.40%  '
$ws.Range('D4').Value = "'0.9999"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = "'304.90"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.69%  '
$ws.Range('D6').Value = "'0.9997"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.07%  '
$ws.Range('D7').Value = "'0.5059"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.29%  '
$ws.Range('D8').Value = "'0.3643"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.89%  '
$ws.Range('D9').Value = "'0.07171"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.47%  '
$ws.Range('D10').Value = "'0.8951"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.85%  '
$ws.Range('D11').Value = "'20.81"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.82%  '
$ws.Range('D12').Value = "'0.07488"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.73%  '
$ws.Range('D13').Value = "'1.836.64"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.98%  '
$ws.Range('D14').Value = "'92.60"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.87%  '
$ws.Range('D15').Value = "'5.232"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.65%  '
$ws.Range('D16').Value = "'1.000"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.02%  '
$ws.Range('D17').Value = "'0.000008488"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.18%  '
$ws.Range('D18').Value = "'14.17"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.45%  '
$ws.Range('D19').Value = "'0.9995"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.05%  '
$ws.Range('D20').Value = "'26.952.60"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.82%  '
$ws.Range('D21').Value = "'5.037"
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Value = "'2.065.21"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.34%  '
$ws.Range('D23').Value = "'10.38"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.67%  '
$ws.Range('D24').Value = "'6.396"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.38%  '
$ws.Range('D25').Value = "'147.42"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.63%  '
$ws.Range('D26').Value = "'1.791"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.61%  '
$ws.Range('D27').Value = "'17.90"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.32%  '
$ws.Range('D28').Value = "'2.074"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.12%  '
$ws.Range('D29').Value = "'113.26"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.43%  '
$ws.Range('D30').Value = "'4.706"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.18%  '
$ws.Range('D31').Value = "'4.682"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.03%  '
$ws.Range('D32').Value = "'0.09250"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.53%  '
$ws.Range('D33').Value = "'0.05103"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.60%  '
$ws.Range('D34').Value = "'0.7550"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.80%  '
$ws.Range('D35').Value = "'2.980"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.41%  '
$ws.Range('D36').Value = "'1.151"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.80%  '
$ws.Range('D37').Value = "'3.271"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +6.83%  '
$ws.Range('D38').Value = "'2.539"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.34%  '
$ws.Range('D39').Value = "'0.01997"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.70%  '
$ws.Range('D40').Value = "'0.5541"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.62%  '
$ws.Range('D41').Value = "'1.072"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.66%  '
$ws.Range('D42').Value = "'118.72"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.18%  '
$ws.Range('D43').Value = "'6.513"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.08%  '
$ws.Range('D44').Value = "'8.521"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.23%  '
$ws.Range('D45').Value = "'0.1472"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.20%  '
$ws.Range('E46').Value = '  +1.20%  '
$ws.Range('D47').Value = "'0.9991"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.07%  '
$ws.Range('D48').Value = "'10.10"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.13%  '
$ws.Range('D49').Value = "'1.563"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.44%  '
$ws.Range('D50').Value = "'36.86"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.05%  '
$ws.Range('E51').Value = '  -2.27%  '
